# Dados_WEGE3.xlsx — "Add files via upload"
#
# VAR_RLO (sheet4) gets two new columns:
#   C = "media geometrica"    -> 5-year rolling GEOMEAN of column B
#   D = "taxa de crescimento" -> period-over-period growth rate of column C
# The VAR_RLO tab also becomes the active/selected sheet (it was BP_WEGE3
# before), with the selection left on D4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VAR_RLO")

# --- Headers -----------------------------------------------------------
$ws.Range("C1").Value = "media geometrica"
$ws.Range("D1").Value = "taxa de crescimento"

# --- Placeholder zeros for the rows before the rolling window is full --
$ws.Range("C2:C5").Value = 0
$ws.Range("D2:D6").Value = 0

# --- Rolling 5-year geometric mean of Receita liquida operacional (B) --
$ws.Range("C6").Formula = "=GEOMEAN(B2:B6)"
$ws.Range("C7:C11").Formula = "=GEOMEAN(B3:B7)"

# --- Year-over-year growth rate of the geometric mean (D) --------------
$ws.Range("D7").Formula = "=C7/C6-1"
$ws.Range("D8:D11").Formula = "=C8/C7-1"

# --- Column widths for the new columns ----------------------------------
$ws.Columns.Item(3).ColumnWidth = 16.5
$ws.Columns.Item(4).ColumnWidth = 18.166666666666668

# --- Make VAR_RLO the active sheet/tab, with D4 selected ----------------
$ws.Activate()
$ws.Range("D4").Select()
